# Updating filtered feeds from workflow
# Appends one new row (row 39) to the "Filtered Feeds" sheet for a newly
# scraped item: a 360dx.com link, its "CDx" keyword, and its title.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

$link = "https://www.360dx.com/cancer/agilent-secures-ivdr-class-c-certification-colorectal-cancer-cdx"
$keywords = "CDx"
$title = "Agilent Secures IVDR Class C Certification for Colorectal Cancer CDx"

$cellA = $ws.Cells.Item($newRow, 1)
$cellB = $ws.Cells.Item($newRow, 2)
$cellC = $ws.Cells.Item($newRow, 3)

$cellA.Value2 = $link
$cellB.Value2 = $keywords
$cellC.Value2 = $title

# Add the external hyperlink for the link cell, then (re)apply the
# worksheet's "Hyperlink" cell style so it visually/structurally matches
# the other link cells in column A.
$ws.Hyperlinks.Add($cellA, $link) | Out-Null
$cellA.Style = "Hyperlink"

Write-Host "Added row $newRow : $link"
